$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the post row for 「あす日が昇るだろう。私は朝が大好きだ」 (row 502),
# shifting all subsequent rows up by one.
$ws.Rows("502").Delete()
